# Generate Report for Archive
# The report rows for "a49fa91a-90f4-44fd-87cc-754bc78e0acf" and
# "cf099737-ea7a-4354-ab2d-e0d79982a66d" (previously rows 5 and 4,
# respectively) are re-sorted by their Latest Handoff Datetime, so the
# cf099737 record now appears first (row 4) and the a49fa91a record
# second (row 5). The cf099737 record's status is also refreshed from
# "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (row 4 <-> row 5 content swap, with cf099737 status
# updated from "Ready for handoff" to "In Translation")
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
$wsOverview.Range("B4").Value = "e2e\cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-09-07 05:23:08"

$wsOverview.Range("A5").Value = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
$wsOverview.Range("B5").Value = "e2e\a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
$wsOverview.Range("G5").Value = "2016-09-07 05:24:11"

# The hyperlinks in column B stay anchored to the same cells (B4/B5) but
# their visible text must follow the swapped content.
$overviewLinks = @()
foreach ($hl in $wsOverview.Hyperlinks) {
    $overviewLinks += $hl
}
foreach ($hl in $overviewLinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "$B$4") {
        $hl.TextToDisplay = "e2e\cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
    } elseif ($addr -eq "$B$5") {
        $hl.TextToDisplay = "e2e\a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn" (row 4 <-> row 5 content swap)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
$wsZhCn.Range("C4").Value = "In Translation"
$wsZhCn.Range("G4").Value = "cf099737-ea7a-4354-ab2d-e0d79982a66d.c0ad95c22bd5146a6597615f7e2e79e1c30ae578.zh-cn.xlf"
$wsZhCn.Range("H4").Value = "2016-09-07 05:22:52"

$wsZhCn.Range("A5").Value = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
$wsZhCn.Range("G5").Value = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.a53e508de1df20ee91cd5532f8ba853bcace8880.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-07 05:23:53"

$zhCnLinks = @()
foreach ($hl in $wsZhCn.Hyperlinks) {
    $zhCnLinks += $hl
}
foreach ($hl in $zhCnLinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "$A$4") {
        $hl.TextToDisplay = "cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
    } elseif ($addr -eq "$A$5") {
        $hl.TextToDisplay = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
    }
}

# ---------------------------------------------------------------------
# Sheet "de-de" (row 4 <-> row 5 content swap)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
$wsDeDe.Range("C4").Value = "In Translation"
$wsDeDe.Range("G4").Value = "cf099737-ea7a-4354-ab2d-e0d79982a66d.c0ad95c22bd5146a6597615f7e2e79e1c30ae578.de-de.xlf"
$wsDeDe.Range("H4").Value = "2016-09-07 05:23:08"

$wsDeDe.Range("A5").Value = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
$wsDeDe.Range("G5").Value = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.a53e508de1df20ee91cd5532f8ba853bcace8880.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-07 05:24:11"

$deDeLinks = @()
foreach ($hl in $wsDeDe.Hyperlinks) {
    $deDeLinks += $hl
}
foreach ($hl in $deDeLinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq "$A$4") {
        $hl.TextToDisplay = "cf099737-ea7a-4354-ab2d-e0d79982a66d.md"
    } elseif ($addr -eq "$A$5") {
        $hl.TextToDisplay = "a49fa91a-90f4-44fd-87cc-754bc78e0acf.md"
    }
}

$wb.Save()
